$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:48 AM"
$ws.Range("G8").Value = "BM"
$ws.Range("C10").Value = "07/21/2025 to 07/27/25"
$ws.Range("G13").Value = "#NO MATCH-3"
